# Update cryptocurrency price/volume/hour data per daily symbol list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'306.49"
$ws.Range("E2").Value = "'-4.08%"
$ws.Range("G2").Value = "'5"
$ws.Range("D3").Value = "'39.84"
$ws.Range("E3").Value = "'-6.87%"
$ws.Range("G3").Value = "'5"
$ws.Range("D4").Value = "'5.048"
$ws.Range("E4").Value = "'-2.46%"
$ws.Range("G4").Value = "'5"
$ws.Range("D5").Value = "'0.07671"
$ws.Range("E5").Value = "'-5.77%"
$ws.Range("G5").Value = "'5"
$ws.Range("D6").Value = "'4.252"
$ws.Range("E6").Value = "'-2.03%"
$ws.Range("G6").Value = "'5"
$ws.Range("D7").Value = "'1.608"
$ws.Range("E7").Value = "'-8.47%"
$ws.Range("G7").Value = "'5"
$ws.Range("D8").Value = "'0.8832"
$ws.Range("E8").Value = "'-7.20%"
$ws.Range("G8").Value = "'5"
$ws.Range("D9").Value = "'0.09978"
$ws.Range("E9").Value = "'-10.13%"
$ws.Range("G9").Value = "'5"
$ws.Range("E10").Value = "'-6.42%"
$ws.Range("G10").Value = "'5"
$ws.Range("D11").Value = "'0.08925"
$ws.Range("E11").Value = "'-4.52%"
$ws.Range("G11").Value = "'5"
$ws.Range("D12").Value = "'0.04392"
$ws.Range("E12").Value = "'-6.17%"
$ws.Range("G12").Value = "'5"
$ws.Range("D13").Value = "'0.1052"
$ws.Range("E13").Value = "'-0.71%"
$ws.Range("G13").Value = "'5"
$ws.Range("D14").Value = "'0.001238"
$ws.Range("E14").Value = "'-3.64%"
$ws.Range("G14").Value = "'5"
$ws.Range("D15").Value = "'0.005885"
$ws.Range("E15").Value = "'-1.81%"
$ws.Range("G15").Value = "'5"
$ws.Range("D16").Value = "'3.358"
$ws.Range("E16").Value = "'-0.24%"
$ws.Range("G16").Value = "'5"
$ws.Range("D17").Value = "'2.462"
$ws.Range("E17").Value = "'-2.99%"
$ws.Range("G17").Value = "'5"
$ws.Range("G18").Value = "'5"
$ws.Range("D19").Value = "'6.991"
$ws.Range("E19").Value = "'-5.71%"
$ws.Range("G19").Value = "'5"
$ws.Range("D20").Value = "'0.1354"
$ws.Range("E20").Value = "'-1.96%"
$ws.Range("G20").Value = "'5"
$ws.Range("D21").Value = "'0.3200"
$ws.Range("E21").Value = "'21.77%"
$ws.Range("G21").Value = "'5"
$ws.Range("D22").Value = "'0.04217"
$ws.Range("E22").Value = "'0.73%"
$ws.Range("G22").Value = "'5"
$ws.Range("D23").Value = "'0.001199"
$ws.Range("E23").Value = "'-4.29%"
$ws.Range("G23").Value = "'5"
$ws.Range("D24").Value = "'0.004056"
$ws.Range("E24").Value = "'-6.01%"
$ws.Range("G24").Value = "'5"
$ws.Range("D25").Value = "'0.0001220"
$ws.Range("E25").Value = "'9.43%"
$ws.Range("G25").Value = "'5"
$ws.Range("E26").Value = "'-0.32%"
$ws.Range("G26").Value = "'5"
$ws.Range("G27").Value = "'5"
$ws.Range("G28").Value = "'5"
$ws.Range("G29").Value = "'5"
$ws.Range("G30").Value = "'5"
$ws.Range("G31").Value = "'5"
$ws.Range("G32").Value = "'5"
$ws.Range("G33").Value = "'5"
$ws.Range("G34").Value = "'5"
$ws.Range("G35").Value = "'5"
$ws.Range("G36").Value = "'5"
$ws.Range("G37").Value = "'5"
$ws.Range("D38").Value = "'0.02341"
$ws.Range("E38").Value = "'-9.49%"
$ws.Range("G38").Value = "'5"
$ws.Range("D39").Value = "'0.05135"
$ws.Range("E39").Value = "'-7.18%"
$ws.Range("G39").Value = "'5"
$ws.Range("D40").Value = "'0.007947"
$ws.Range("E40").Value = "'1.65%"
$ws.Range("G40").Value = "'5"
$ws.Range("D41").Value = "'0.1318"
$ws.Range("E41").Value = "'-5.36%"
$ws.Range("G41").Value = "'5"
$ws.Range("D42").Value = "'0.006663"
$ws.Range("E42").Value = "'0.64%"
$ws.Range("G42").Value = "'5"
$ws.Range("D43").Value = "'0.001985"
$ws.Range("E43").Value = "'-6.49%"
$ws.Range("G43").Value = "'5"
$ws.Range("D44").Value = "'0.008471"
$ws.Range("E44").Value = "'0.01%"
$ws.Range("G44").Value = "'5"
$ws.Range("D45").Value = "'0.3040"
$ws.Range("E45").Value = "'-12.06%"
$ws.Range("G45").Value = "'5"
$ws.Range("D46").Value = "'0.00006535"
$ws.Range("E46").Value = "'-6.55%"
$ws.Range("G46").Value = "'5"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("E47").Value = "'-0.40%"
$ws.Range("G47").Value = "'5"
$ws.Range("B48").Value = "'CoinbaseStockToken"
$ws.Range("C48").Value = "'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").Value = "'0.007000"
$ws.Range("E48").Value = "'97.50%"
$ws.Range("G48").Value = "'5"
$ws.Range("B49").Value = "'BOLO"
$ws.Range("C49").Value = "'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").Value = "'0.003392"
$ws.Range("E49").Value = "'-2.43%"
$ws.Range("G49").Value = "'5"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("E50").Value = "'-0.40%"
$ws.Range("G50").Value = "'5"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("E51").Value = "'-0.40%"
$ws.Range("G51").Value = "'5"
